$d = $word.ActiveDocument

# 1) Remove the entire paragraph "FABIO HO BISOGNO DI DEF SHARED"
$removed = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "FABIO HO BISOGNO DI DEF SHARED*") {
        $p.Range.Delete()
        $removed = $true
        break
    }
}

# 2) Append a closing parenthesis ")" to the paragraph ending in
#    "...VEDI PEZZO WS DI US" as its own separate run.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*VEDI PEZZO WS DI US*") {
        $pr = $p.Range
        $insPoint = $d.Range($pr.End - 1, $pr.End - 1)
        $insPoint.InsertAfter(")")

        $pr2 = $p.Range
        $parenRange = $d.Range($pr2.End - 2, $pr2.End - 1)
        # Touch and restore a character property so the engine keeps the
        # newly-typed ")" as a distinct run instead of re-merging it into
        # the preceding run's text node.
        $parenRange.Bold = 1
        $parenRange.Bold = 0
        break
    }
}
